$wb = $excel.ActiveWorkbook

# --- Sheet "UISA" (sheet1): correct B58 and append new weekly row 59 ---
$ws1 = $wb.Worksheets.Item("UISA")
$ws1.Range("B58").Value = 769000

# Copy formatting from the row above (dates use style 4, values use style 3)
# so the new row inherits the same look without minting new styles.
$ws1.Range("A58:B58").Copy()
$ws1.Range("A59:B59").PasteSpecial(-4122)
$ws1.Range("A59").Value = 44296
$ws1.Range("B59").Value = 576000

# --- Sheet "ContClaims" (sheet2): correct B57 and append new weekly row 58 ---
$ws2 = $wb.Worksheets.Item("ContClaims")
$ws2.Range("B57").Value = 3727000

$ws2.Range("A57:B57").Copy()
$ws2.Range("A58:B58").PasteSpecial(-4122)
$ws2.Range("A58").Value = 44289
$ws2.Range("B58").Value = 3731000
